$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 62; everything at/below row 62 shifts down by one
# (old row 62 -> 63, ..., old row 94 -> 95).
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new record.
$ws.Cells.Item(62, 1).Value = 6
$ws.Cells.Item(62, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(62, 3).Value = "Metropolitana"
$ws.Cells.Item(62, 4).Value = 44452
$ws.Cells.Item(62, 5).Value = 13
$ws.Cells.Item(62, 6).Value = 100112001
$ws.Cells.Item(62, 7).Value = "Berenjena"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 10).Value = 730
$ws.Cells.Item(62, 11).Value = 6000
$ws.Cells.Item(62, 12).Value = 7000
$ws.Cells.Item(62, 13).Value = 6479
$ws.Cells.Item(62, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(62, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(62, 16).Value = 130
$ws.Cells.Item(62, 17).Value = 50
$ws.Cells.Item(62, 18).Value = "Hortaliza"
